# Update gh-pages to output generated at 456a3b4
# Applies updated view/follower style counters to the 展览 (sheet1),
# 演出 (sheet2) and 全部类型 (sheet4) worksheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1429
$ws.Range("G3").Value = 68
$ws.Range("F4").Value = 20103
$ws.Range("F5").Value = 799
$ws.Range("F9").Value = 7581
$ws.Range("F15").Value = 120
$ws.Range("F16").Value = 12
$ws.Range("F20").Value = 423
$ws.Range("F24").Value = 71
$ws.Range("F31").Value = 5219
$ws.Range("F34").Value = 2851
$ws.Range("F38").Value = 12647
$ws.Range("F40").Value = 83
$ws.Range("F42").Value = 56
$ws.Range("F43").Value = 269
$ws.Range("F44").Value = 368

# ---- Sheet "演出" ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 179

# ---- Sheet "全部类型" ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1429
$ws.Range("G3").Value = 68
$ws.Range("F4").Value = 20103
$ws.Range("F5").Value = 799
$ws.Range("F9").Value = 7581
$ws.Range("F15").Value = 120
$ws.Range("F16").Value = 12
$ws.Range("F20").Value = 423
$ws.Range("F24").Value = 71
$ws.Range("F31").Value = 179
$ws.Range("F36").Value = 2851
$ws.Range("F40").Value = 12647
$ws.Range("F42").Value = 83
$ws.Range("F44").Value = 56
$ws.Range("F45").Value = 269
$ws.Range("F46").Value = 368
